$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.337.28'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '3.944.85'
$ws.Range("E3").Value = '  +4.50%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '471.42'
$ws.Range("E5").Value = '  +8.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.59'
$ws.Range("E6").Value = '  +5.36%  '
$ws.Range("E7").Value = '  +1.31%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.736'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("E10").Value = '  +9.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000341'
$ws.Range("E11").Value = '  +9.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.57'
$ws.Range("E12").Value = '  +1.55%  '
$ws.Range("D13").Value = '4.569.47'
$ws.Range("E13").Value = '  +4.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.45'
$ws.Range("E14").Value = '  +0.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.15'
$ws.Range("E15").Value = '  +2.43%  '
$ws.Range("D16").Value = '3.925.91'
$ws.Range("E16").Value = '  +4.01%  '
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.91'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("E19").Value = '  +2.25%  '
$ws.Range("D20").Value = '67.563.60'
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.50'
$ws.Range("E21").Value = '  +7.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.41'
$ws.Range("E22").Value = '  +5.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.53'
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.68'
$ws.Range("E24").Value = '  +2.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.62'
$ws.Range("E25").Value = '  +7.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.80'
$ws.Range("E26").Value = '  +5.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.32'
$ws.Range("E27").Value = '  +5.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.84'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '720.36'
$ws.Range("E29").Value = '  -0.69%  '
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.81'
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '42.77'
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.98'
$ws.Range("E34").Value = '  +3.28%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.151'
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").Value = '0.0₃0808'
$ws.Range("E36").Value = '  +21.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.37'
$ws.Range("E38").Value = '  -4.56%  '
$ws.Range("E39").Value = '  +0.76%  '
$ws.Range("E40").Value = '  +2.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.56'
$ws.Range("E41").Value = '  +7.69%  '
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.340'
$ws.Range("E43").Value = '  +3.20%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.58'
$ws.Range("E44").Value = '  -6.38%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("E46").Value = '  +4.99%  '
$ws.Range("E47").Value = '  +6.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '147.56'
$ws.Range("E48").Value = '  +4.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.18'
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("E50").Value = '  +2.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.92'
$ws.Range("E51").Value = '  +5.11%  '
